$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per diff (row number corresponds to sheet row)
$ws.Range("C2").Value = 12.4

$ws.Range("B3").Value = -1.3
$ws.Range("C3").Value = 7.1

$ws.Range("C4").Value = 6.8

$ws.Range("C5").Value = 11

$ws.Range("C6").Value = 11.4

$ws.Range("C7").Value = 12.6

$ws.Range("C8").Value = 5.1

$ws.Range("C10").Value = 6.6

$ws.Range("C11").Value = 12

$ws.Range("C18").Value = 7.8

$ws.Range("C20").Value = 6.2

$ws.Range("C22").Value = 11

$ws.Range("C24").Value = -0
